# Auto-generated edit script: update cryptos price/volume data per commit
# "Updated cryptos list on Mon Mar 20 09:15:03 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.172.88"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "1.782.54"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'336.46"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'0.9951"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "'0.3821"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.3428"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'47.37"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'1.155"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").Value = "'0.07405"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'23.12"
$ws.Range("E12").Value = "  +7.15%  "
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'6.419"
$ws.Range("D15").Value = "1.783.56"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'7.147"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'0.00001081"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "'0.06655"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'82.72"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'0.9973"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'17.47"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").Value = "'6.425"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "28.249.57"
$ws.Range("E23").Value = "  +4.19%  "
$ws.Range("D24").Value = "'12.08"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'2.367"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'20.87"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "'1.426"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").Value = "'2.407"
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("D29").Value = "'154.50"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "1.987.61"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "'135.40"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'6.156"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("D33").Value = "'3.968"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'0.08758"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "'12.77"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "'0.02422"
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").Value = "'0.6865"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("D38").Value = "'5.332"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "'0.06355"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "'0.2176"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.505"
$ws.Range("E41").Value = "  -7.37%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.240"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'8.309"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").Value = "'14.38"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").Value = "'0.9964"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "'0.6309"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "'3.852"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'132.20"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "'2.090"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'0.07500"
$ws.Range("E50").Value = "  +5.81%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.205"
$ws.Range("E51").Value = "  +8.59%  "
